$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the shared-string header in B1: append the "(1416188)" suffix
$ws.Range("B1").Value = "consumer_inflation_expectation (1416188)"

# 2) Widen column B from 30.0 to 36.5 (character width units, accounting for
#    Excel's internal padding offset so the saved <col width="..."/> matches)
$ws.Columns("B").ColumnWidth = 35.666666666666664

# 3) Append two new monthly data rows (218, 219) after the existing last row (217)
$ws.Range("A218").Value = 45170
$ws.Range("B218").Value = 7.3
$ws.Range("A219").Value = 45200
$ws.Range("B219").Value = 7.5

# Copy the number formatting/style from the preceding data row (217) onto the
# two new rows so they carry the same date/number styles as the rest of the column
$ws.Range("A217:B217").Copy()
$ws.Range("A218:B219").PasteSpecial(-4122)
